$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "2-2090-O"
$ws.Range("E1").Value = "HTML-2030-P"
$ws.Range("F1").Value = "777-2080-P"
$ws.Range("G1").Value = "PHP"
$ws.Range("H1").Value = "88-2015"

$ws.Range("D2").Select()
